$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet lists weekly Achicoria price records (rows 2-21). This commit
# reshuffles the Fecha/Volumen/Precio/Origen fields among the existing rows
# (the row order of these fields no longer matches the previous date order).
# For each target row below, the values are exactly the ones that used to
# live in the "original row" noted in the comment (taken from the workbook
# before this edit).

# Row 2 <- original Row 8
$ws.Range("D2").Value2  = 44232
$ws.Range("J2").Value2  = 250
$ws.Range("K2").Value2  = 5000
$ws.Range("L2").Value2  = 6000
$ws.Range("M2").Value2  = 5500
$ws.Range("O2").Value2  = "Provincia de Quillota"
$ws.Range("P2").Value2  = 344

# Row 3 <- original Row 2
$ws.Range("D3").Value2  = 44230
$ws.Range("J3").Value2  = 250
$ws.Range("K3").Value2  = 5000
$ws.Range("L3").Value2  = 6000
$ws.Range("M3").Value2  = 5500
$ws.Range("O3").Value2  = "Provincia de Quillota"
$ws.Range("P3").Value2  = 344

# Row 4 <- original Row 12
$ws.Range("D4").Value2  = 44855
$ws.Range("J4").Value2  = 70
$ws.Range("K4").Value2  = 6000
$ws.Range("L4").Value2  = 7000
$ws.Range("M4").Value2  = 6500
$ws.Range("O4").Value2  = "Provincia de Quillota"
$ws.Range("P4").Value2  = 406

# Row 5 <- original Row 19
$ws.Range("D5").Value2  = 44231
$ws.Range("J5").Value2  = 250
$ws.Range("K5").Value2  = 5000
$ws.Range("L5").Value2  = 6000
$ws.Range("M5").Value2  = 5500
$ws.Range("O5").Value2  = "Provincia de Quillota"
$ws.Range("P5").Value2  = 344

# Row 6 <- original Row 16
$ws.Range("D6").Value2  = 44883
$ws.Range("J6").Value2  = 180
$ws.Range("K6").Value2  = 7000
$ws.Range("L6").Value2  = 8000
$ws.Range("M6").Value2  = 7500
$ws.Range("O6").Value2  = "Provincia de Quillota"
$ws.Range("P6").Value2  = 469

# Row 7 <- original Row 9
$ws.Range("D7").Value2  = 44873
$ws.Range("J7").Value2  = 250
$ws.Range("K7").Value2  = 8000
$ws.Range("L7").Value2  = 8000
$ws.Range("M7").Value2  = 8000
$ws.Range("O7").Value2  = "Provincia de Quillota"
$ws.Range("P7").Value2  = 500

# Row 8 <- original Row 3
$ws.Range("D8").Value2  = 44208
$ws.Range("J8").Value2  = 160
$ws.Range("K8").Value2  = 5000
$ws.Range("L8").Value2  = 6000
$ws.Range("M8").Value2  = 5500
$ws.Range("O8").Value2  = "Provincia de Quillota"
$ws.Range("P8").Value2  = 344

# Row 9 <- original Row 15
$ws.Range("D9").Value2  = 44186
$ws.Range("J9").Value2  = 160
$ws.Range("K9").Value2  = 5000
$ws.Range("L9").Value2  = 6000
$ws.Range("M9").Value2  = 5500
$ws.Range("O9").Value2  = "Provincia de Quillota"
$ws.Range("P9").Value2  = 344

# Row 10 <- original Row 17
$ws.Range("D10").Value2 = 44189
$ws.Range("J10").Value2 = 250
$ws.Range("K10").Value2 = 5000
$ws.Range("L10").Value2 = 6000
$ws.Range("M10").Value2 = 5500
$ws.Range("O10").Value2 = "Provincia de Quillota"
$ws.Range("P10").Value2 = 344

# Row 11 <- original Row 21
$ws.Range("D11").Value2 = 44846
$ws.Range("J11").Value2 = 250
$ws.Range("K11").Value2 = 5000
$ws.Range("L11").Value2 = 5000
$ws.Range("M11").Value2 = 5000
$ws.Range("O11").Value2 = "Provincia de Quillota"
$ws.Range("P11").Value2 = 312

# Row 12 <- original Row 6
$ws.Range("D12").Value2 = 44882
$ws.Range("J12").Value2 = 70
$ws.Range("K12").Value2 = 7000
$ws.Range("L12").Value2 = 7000
$ws.Range("M12").Value2 = 7000
$ws.Range("O12").Value2 = "Provincia de Quillota"
$ws.Range("P12").Value2 = 438

# Row 13 <- original Row 18
$ws.Range("D13").Value2 = 44187
$ws.Range("J13").Value2 = 160
$ws.Range("K13").Value2 = 5000
$ws.Range("L13").Value2 = 6000
$ws.Range("M13").Value2 = 5500
$ws.Range("O13").Value2 = "Provincia de Quillota"
$ws.Range("P13").Value2 = 344

# Row 14 <- original Row 11
$ws.Range("D14").Value2 = 44204
$ws.Range("J14").Value2 = 430
$ws.Range("K14").Value2 = 5000
$ws.Range("L14").Value2 = 6000
$ws.Range("M14").Value2 = 5500
$ws.Range("O14").Value2 = "Provincia de Quillota"
$ws.Range("P14").Value2 = 344

# Row 15 <- original Row 13
$ws.Range("D15").Value2 = 44251
$ws.Range("J15").Value2 = 120
$ws.Range("K15").Value2 = 5000
$ws.Range("L15").Value2 = 5000
$ws.Range("M15").Value2 = 5000
$ws.Range("O15").Value2 = "Región Metropolitana"
$ws.Range("P15").Value2 = 312

# Row 16 <- original Row 7
$ws.Range("D16").Value2 = 44210
$ws.Range("J16").Value2 = 340
$ws.Range("K16").Value2 = 5000
$ws.Range("L16").Value2 = 6000
$ws.Range("M16").Value2 = 5500
$ws.Range("O16").Value2 = "Provincia de Quillota"
$ws.Range("P16").Value2 = 344

# Row 17 <- original Row 5
$ws.Range("D17").Value2 = 44875
$ws.Range("J17").Value2 = 90
$ws.Range("K17").Value2 = 7000
$ws.Range("L17").Value2 = 7000
$ws.Range("M17").Value2 = 7000
$ws.Range("O17").Value2 = "Provincia de Quillota"
$ws.Range("P17").Value2 = 438

# Row 18 <- original Row 4
$ws.Range("D18").Value2 = 44236
$ws.Range("J18").Value2 = 180
$ws.Range("K18").Value2 = 4000
$ws.Range("L18").Value2 = 4500
$ws.Range("M18").Value2 = 4167
$ws.Range("O18").Value2 = "Región Metropolitana"
$ws.Range("P18").Value2 = 260

# Row 19 <- original Row 20
$ws.Range("D19").Value2 = 44292
$ws.Range("J19").Value2 = 90
$ws.Range("K19").Value2 = 6000
$ws.Range("L19").Value2 = 6000
$ws.Range("M19").Value2 = 6000
$ws.Range("O19").Value2 = "Región Metropolitana"
$ws.Range("P19").Value2 = 375

# Row 20 <- original Row 10
$ws.Range("D20").Value2 = 44188
$ws.Range("J20").Value2 = 210
$ws.Range("K20").Value2 = 5000
$ws.Range("L20").Value2 = 6000
$ws.Range("M20").Value2 = 5500
$ws.Range("O20").Value2 = "Provincia de Quillota"
$ws.Range("P20").Value2 = 344

# Row 21 <- original Row 14
$ws.Range("D21").Value2 = 44215
$ws.Range("J21").Value2 = 250
$ws.Range("K21").Value2 = 5000
$ws.Range("L21").Value2 = 6000
$ws.Range("M21").Value2 = 5500
$ws.Range("O21").Value2 = "Provincia de Quillota"
$ws.Range("P21").Value2 = 344
